# Updates crypto prices / 1h volume percentages (and restores the original
# ranking order for rows 16/17 = WrappedEther/ShibaInu and rows 50/51 =
# Cronos/EOS), matching the refreshed data pulled by the scheduled
# GitHub Actions job.
#
# Column D holds price text that sometimes *looks* numeric (e.g. "0.3680",
# "0.00001137"); a plain `.Value = "..."` assignment lets Excel auto-detect
# those as real numbers, which silently drops trailing zeros / switches to
# scientific notation and changes the cell's stored type from text to
# number. To avoid that, each D-column write temporarily forces the cell to
# Text format, assigns the literal string, then restores the cell's
# (default) style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.386.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.570.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.92%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3680"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.52"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3376"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07557"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.059"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.852"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.21%  "

# Row 16/17: swap WrappedEther <-> ShibaInu back to the canonical order
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001137"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.25%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.573.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06676"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.38%  "

$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.216"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.394.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.403"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.949"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "145.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.921"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.749.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.238"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.968"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9840"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -12.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08439"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02539"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2294"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06483"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.469"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.14%  "

$ws.Range("E41").Value = "  -12.26%  "

$ws.Range("E42").Value = "  -5.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6368"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.90%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5983"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.768"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.102"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.14%  "

# Row 50/51: swap EOS <-> Cronos back to the canonical order
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07265"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.74%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.188"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.57%  "
